# Update "想去人数" (column F) values on the 展览 and 全部类型 sheets
# to reflect newly scraped counts, per commit:
#   "Update gh-pages to output generated at 456a3b4"

$wb = $excel.ActiveWorkbook

# Row -> new F-column value for 展览 (sheet 1) and 全部类型 (sheet 4).
# Both sheets share the same underlying rows except the last one (F23),
# which diverges by one between the two sheets in this refresh.
$commonUpdates = @{
    2  = 351
    4  = 10624
    6  = 966
    7  = 113
    8  = 1313
    9  = 8228
    10 = 32
    15 = 3266
    18 = 749
    20 = 1056
    22 = 99
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $commonUpdates.Keys) {
        $ws.Range("F$row").Value = $commonUpdates[$row]
    }
}

# F23 differs between the two sheets in the refreshed data.
$wb.Worksheets.Item("展览").Range("F23").Value = 1711
$wb.Worksheets.Item("全部类型").Range("F23").Value = 1712
